# adding Tp 28-30 from OR
# Append a new data row (row 81) to Sheet1, mirroring the pattern of the
# existing rows: a date in column A, CRM/Batch values in B/C, the
# "% off" formula in D, the constant 169 in E, and a note in F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (date number format/style) of the row above down
# into the new row so A81 keeps the same date style as A80.
$ws.Range("A80").Copy() | Out-Null
$ws.Range("A81").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New sample row for 12/17/2019.
$ws.Range("A81").Value = 43816
$ws.Range("B81").Value = 2219.6091874158001
$ws.Range("C81").Value = 2207.0300000000002
$ws.Range("D81").Formula = "=100*(B81-C81)/C81"
$ws.Range("E81").Value = 169
$ws.Range("F81").Value = "New CRM opened 12/11/2023"

# Match the author's final selection in the sheet view.
$ws.Range("F80:F81").Select()
